$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: confusion_matrix (update text values) ---
$ws.Range("B2").Value = "[[12242    97]`n [  243   310]]"
$ws.Range("C2").Value = "[[12292    47]`n [  319   234]]"
$ws.Range("D2").Value = "[[12285    54]`n [  321   232]]"
$ws.Range("E2").Value = "[[12162   177]`n [  280   273]]"
$ws.Range("F2").Value = "[[12302    37]`n [  425   128]]"
$ws.Range("G2").Value = "[[12300    39]`n [  381   172]]"
# Re-run autofit on row 2 so the multi-line text doesn't leave a stray
# custom row height behind (matches the original unmodified row height).
$ws.Rows("2").AutoFit()

# --- Row 3: accuracy_score (update numeric values) ---
$ws.Range("B3").Value = 0.9736270555383183
$ws.Range("C3").Value = 0.9716103009618368
$ws.Range("D3").Value = 0.9709121936084394
$ws.Range("E3").Value = 0.9645516599441514
$ws.Range("F3").Value = 0.9641638225255973
$ws.Range("G3").Value = 0.9674216568414521

# --- Row 4: f1_score (update numeric values) ---
$ws.Range("B4").Value = 0.6458333333333334
$ws.Range("C4").Value = 0.5611510791366906
$ws.Range("D4").Value = 0.5530393325387366
$ws.Range("E4").Value = 0.5443668993020938
$ws.Range("F4").Value = 0.3565459610027855
$ws.Range("G4").Value = 0.450261780104712

# --- Row 5: new "time" row ---
# Copy formatting from A4 (bold, bordered, centered/top-aligned header style)
# onto the new A5 cell, then set its text.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "time"

$ws.Range("B5").Value = 7.668647050857544
$ws.Range("C5").Value = 18.13580298423767
$ws.Range("D5").Value = 15.49525141716003
$ws.Range("E5").Value = 7.885470628738403
$ws.Range("F5").Value = 57.8807532787323
$ws.Range("G5").Value = 71.08479619026184
